$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Matriz_Resultados")
$ws.Range("B2").Value2 = [double]"0"
$ws.Range("C2").Value2 = [double]"0"
$ws.Range("D2").Value2 = [double]"0"
$ws.Range("E2").Value2 = [double]"0"
$ws.Range("F2").Value2 = [double]"0"
$ws.Range("G2").Value2 = [double]"0"
$ws.Range("H2").Value2 = [double]"0"
$ws.Range("I2").Value2 = [double]"0"
$ws.Range("J2").Value2 = [double]"0"
$ws.Range("B3").Value2 = [double]"0"
$ws.Range("C3").Value2 = [double]"0"
$ws.Range("D3").Value2 = [double]"1"
$ws.Range("E3").Value2 = [double]"1"
$ws.Range("F3").Value2 = [double]"0"
$ws.Range("G3").Value2 = [double]"0"
$ws.Range("H3").Value2 = [double]"0"
$ws.Range("I3").Value2 = [double]"0"
$ws.Range("J3").Value2 = [double]"1"
$ws.Range("B4").Value2 = [double]"0"
$ws.Range("C4").Value2 = [double]"-1"
$ws.Range("D4").Value2 = [double]"0"
$ws.Range("E4").Value2 = [double]"1"
$ws.Range("F4").Value2 = [double]"0"
$ws.Range("G4").Value2 = [double]"0"
$ws.Range("H4").Value2 = [double]"0"
$ws.Range("I4").Value2 = [double]"0"
$ws.Range("J4").Value2 = [double]"1"
$ws.Range("B5").Value2 = [double]"0"
$ws.Range("C5").Value2 = [double]"-1"
$ws.Range("D5").Value2 = [double]"-1"
$ws.Range("E5").Value2 = [double]"0"
$ws.Range("F5").Value2 = [double]"0"
$ws.Range("G5").Value2 = [double]"0"
$ws.Range("H5").Value2 = [double]"0"
$ws.Range("I5").Value2 = [double]"0"
$ws.Range("J5").Value2 = [double]"0"
$ws.Range("B6").Value2 = [double]"0"
$ws.Range("C6").Value2 = [double]"0"
$ws.Range("D6").Value2 = [double]"0"
$ws.Range("E6").Value2 = [double]"0"
$ws.Range("F6").Value2 = [double]"0"
$ws.Range("G6").Value2 = [double]"0"
$ws.Range("H6").Value2 = [double]"0"
$ws.Range("I6").Value2 = [double]"0"
$ws.Range("J6").Value2 = [double]"0"
$ws.Range("B7").Value2 = [double]"0"
$ws.Range("C7").Value2 = [double]"0"
$ws.Range("D7").Value2 = [double]"0"
$ws.Range("E7").Value2 = [double]"0"
$ws.Range("F7").Value2 = [double]"0"
$ws.Range("G7").Value2 = [double]"0"
$ws.Range("H7").Value2 = [double]"0"
$ws.Range("I7").Value2 = [double]"0"
$ws.Range("J7").Value2 = [double]"1"
$ws.Range("B8").Value2 = [double]"0"
$ws.Range("C8").Value2 = [double]"0"
$ws.Range("D8").Value2 = [double]"0"
$ws.Range("E8").Value2 = [double]"0"
$ws.Range("F8").Value2 = [double]"0"
$ws.Range("G8").Value2 = [double]"0"
$ws.Range("H8").Value2 = [double]"0"
$ws.Range("I8").Value2 = [double]"0"
$ws.Range("J8").Value2 = [double]"1"
$ws.Range("B9").Value2 = [double]"0"
$ws.Range("C9").Value2 = [double]"0"
$ws.Range("D9").Value2 = [double]"0"
$ws.Range("E9").Value2 = [double]"0"
$ws.Range("F9").Value2 = [double]"0"
$ws.Range("G9").Value2 = [double]"0"
$ws.Range("H9").Value2 = [double]"0"
$ws.Range("I9").Value2 = [double]"0"
$ws.Range("J9").Value2 = [double]"1"
$ws.Range("B10").Value2 = [double]"0"
$ws.Range("C10").Value2 = [double]"-1"
$ws.Range("D10").Value2 = [double]"-1"
$ws.Range("E10").Value2 = [double]"0"
$ws.Range("F10").Value2 = [double]"0"
$ws.Range("G10").Value2 = [double]"-1"
$ws.Range("H10").Value2 = [double]"-1"
$ws.Range("I10").Value2 = [double]"-1"
$ws.Range("J10").Value2 = [double]"0"

$ws = $wb.Worksheets.Item("P_valores")
$ws.Range("B2").Value2 = [double]"1"
$ws.Range("C2").Value2 = [double]"0.005028717609074862"
$ws.Range("D2").Value2 = [double]"0.007001104705730343"
$ws.Range("E2").Value2 = [double]"0.02125018373255627"
$ws.Range("F2").Value2 = [double]"0.04943003968366444"
$ws.Range("G2").Value2 = [double]"0.005944189703441571"
$ws.Range("H2").Value2 = [double]"0.005852768740049141"
$ws.Range("I2").Value2 = [double]"0.007308422022619165"
$ws.Range("J2").Value2 = [double]"0.05908350867441703"
$ws.Range("B3").Value2 = [double]"0.005028717609074862"
$ws.Range("C3").Value2 = [double]"1"
$ws.Range("D3").Value2 = [double]"1.034179414638459E-05"
$ws.Range("E3").Value2 = [double]"0.0003475271253374768"
$ws.Range("F3").Value2 = [double]"0.004262915406137324"
$ws.Range("G3").Value2 = [double]"0.005106953696461858"
$ws.Range("H3").Value2 = [double]"0.005726811050846381"
$ws.Range("I3").Value2 = [double]"0.07032043600425908"
$ws.Range("J3").Value2 = [double]"0.0002329964007945495"
$ws.Range("B4").Value2 = [double]"0.007001104705730343"
$ws.Range("C4").Value2 = [double]"1.034179414638459E-05"
$ws.Range("D4").Value2 = [double]"1"
$ws.Range("E4").Value2 = [double]"0.00103288146616709"
$ws.Range("F4").Value2 = [double]"0.006224403040201221"
$ws.Range("G4").Value2 = [double]"0.02237620212689739"
$ws.Range("H4").Value2 = [double]"0.02740555201311223"
$ws.Range("I4").Value2 = [double]"0.5212429128041012"
$ws.Range("J4").Value2 = [double]"0.0003773948068874766"
$ws.Range("B5").Value2 = [double]"0.02125018373255627"
$ws.Range("C5").Value2 = [double]"0.0003475271253374768"
$ws.Range("D5").Value2 = [double]"0.00103288146616709"
$ws.Range("E5").Value2 = [double]"1"
$ws.Range("F5").Value2 = [double]"0.02234541566585713"
$ws.Range("G5").Value2 = [double]"0.1304225603348284"
$ws.Range("H5").Value2 = [double]"0.05504241534761856"
$ws.Range("I5").Value2 = [double]"0.03742379444293831"
$ws.Range("J5").Value2 = [double]"0.002645474255317248"
$ws.Range("B6").Value2 = [double]"0.04943003968366444"
$ws.Range("C6").Value2 = [double]"0.004262915406137324"
$ws.Range("D6").Value2 = [double]"0.006224403040201221"
$ws.Range("E6").Value2 = [double]"0.02234541566585713"
$ws.Range("F6").Value2 = [double]"1"
$ws.Range("G6").Value2 = [double]"0.005085566982925727"
$ws.Range("H6").Value2 = [double]"0.004990042767744951"
$ws.Range("I6").Value2 = [double]"0.006305592455507547"
$ws.Range("J6").Value2 = [double]"0.084771560727958"
$ws.Range("B7").Value2 = [double]"0.005944189703441571"
$ws.Range("C7").Value2 = [double]"0.005106953696461858"
$ws.Range("D7").Value2 = [double]"0.02237620212689739"
$ws.Range("E7").Value2 = [double]"0.1304225603348284"
$ws.Range("F7").Value2 = [double]"0.005085566982925727"
$ws.Range("G7").Value2 = [double]"1"
$ws.Range("H7").Value2 = [double]"0.2033400490351587"
$ws.Range("I7").Value2 = [double]"0.1658016181642061"
$ws.Range("J7").Value2 = [double]"6.416884939630663E-05"
$ws.Range("B8").Value2 = [double]"0.005852768740049141"
$ws.Range("C8").Value2 = [double]"0.005726811050846381"
$ws.Range("D8").Value2 = [double]"0.02740555201311223"
$ws.Range("E8").Value2 = [double]"0.05504241534761856"
$ws.Range("F8").Value2 = [double]"0.004990042767744951"
$ws.Range("G8").Value2 = [double]"0.2033400490351587"
$ws.Range("H8").Value2 = [double]"1"
$ws.Range("I8").Value2 = [double]"0.2342488368558631"
$ws.Range("J8").Value2 = [double]"4.066635277721886E-05"
$ws.Range("B9").Value2 = [double]"0.007308422022619165"
$ws.Range("C9").Value2 = [double]"0.07032043600425908"
$ws.Range("D9").Value2 = [double]"0.5212429128041012"
$ws.Range("E9").Value2 = [double]"0.03742379444293831"
$ws.Range("F9").Value2 = [double]"0.006305592455507547"
$ws.Range("G9").Value2 = [double]"0.1658016181642061"
$ws.Range("H9").Value2 = [double]"0.2342488368558631"
$ws.Range("I9").Value2 = [double]"1"
$ws.Range("J9").Value2 = [double]"0.001303516221026868"
$ws.Range("B10").Value2 = [double]"0.05908350867441703"
$ws.Range("C10").Value2 = [double]"0.0002329964007945495"
$ws.Range("D10").Value2 = [double]"0.0003773948068874766"
$ws.Range("E10").Value2 = [double]"0.002645474255317248"
$ws.Range("F10").Value2 = [double]"0.084771560727958"
$ws.Range("G10").Value2 = [double]"6.416884939630663E-05"
$ws.Range("H10").Value2 = [double]"4.066635277721886E-05"
$ws.Range("I10").Value2 = [double]"0.001303516221026868"
$ws.Range("J10").Value2 = [double]"1"

$ws = $wb.Worksheets.Item("Estadisticos_DM")
$ws.Range("B2").Value2 = [double]"0"
$ws.Range("C2").Value2 = [double]"3.116377831448006"
$ws.Range("D2").Value2 = [double]"2.974159401053307"
$ws.Range("E2").Value2 = [double]"2.480537185164833"
$ws.Range("F2").Value2 = [double]"2.079533862214786"
$ws.Range("G2").Value2 = [double]"3.044716695131391"
$ws.Range("H2").Value2 = [double]"3.051376050698078"
$ws.Range("I2").Value2 = [double]"2.955561379033963"
$ws.Range("J2").Value2 = [double]"1.990661243152962"
$ws.Range("B3").Value2 = [double]"-3.116377831448006"
$ws.Range("C3").Value2 = [double]"0"
$ws.Range("D3").Value2 = [double]"-5.67964545923587"
$ws.Range("E3").Value2 = [double]"-4.225869615018513"
$ws.Range("F3").Value2 = [double]"-3.186761764257116"
$ws.Range("G3").Value2 = [double]"-3.109780918826911"
$ws.Range("H3").Value2 = [double]"-3.060717104569447"
$ws.Range("I3").Value2 = [double]"-1.902179684934643"
$ws.Range("J3").Value2 = [double]"-4.389490135705872"
$ws.Range("B4").Value2 = [double]"-2.974159401053307"
$ws.Range("C4").Value2 = [double]"5.67964545923587"
$ws.Range("D4").Value2 = [double]"0"
$ws.Range("E4").Value2 = [double]"-3.778793017094997"
$ws.Range("F4").Value2 = [double]"-3.024902456076612"
$ws.Range("G4").Value2 = [double]"-2.456782607457124"
$ws.Range("H4").Value2 = [double]"-2.362636269252887"
$ws.Range("I4").Value2 = [double]"-0.6518618050328975"
$ws.Range("J4").Value2 = [double]"-4.192119868523634"
$ws.Range("B5").Value2 = [double]"-2.480537185164833"
$ws.Range("C5").Value2 = [double]"4.225869615018513"
$ws.Range("D5").Value2 = [double]"3.778793017094997"
$ws.Range("E5").Value2 = [double]"0"
$ws.Range("F5").Value2 = [double]"-2.457417157072955"
$ws.Range("G5").Value2 = [double]"1.571139172864607"
$ws.Range("H5").Value2 = [double]"2.02616217426285"
$ws.Range("I5").Value2 = [double]"2.214968684321045"
$ws.Range("J5").Value2 = [double]"-3.388081490838152"
$ws.Range("B6").Value2 = [double]"-2.079533862214786"
$ws.Range("C6").Value2 = [double]"3.186761764257116"
$ws.Range("D6").Value2 = [double]"3.024902456076612"
$ws.Range("E6").Value2 = [double]"2.457417157072955"
$ws.Range("F6").Value2 = [double]"0"
$ws.Range("G6").Value2 = [double]"3.111574516354902"
$ws.Range("H6").Value2 = [double]"3.119675605408095"
$ws.Range("I6").Value2 = [double]"3.01932165680969"
$ws.Range("J6").Value2 = [double]"1.805059303861617"
$ws.Range("B7").Value2 = [double]"-3.044716695131391"
$ws.Range("C7").Value2 = [double]"3.109780918826911"
$ws.Range("D7").Value2 = [double]"2.456782607457124"
$ws.Range("E7").Value2 = [double]"-1.571139172864607"
$ws.Range("F7").Value2 = [double]"-3.111574516354902"
$ws.Range("G7").Value2 = [double]"0"
$ws.Range("H7").Value2 = [double]"1.311114057127361"
$ws.Range("I7").Value2 = [double]"1.433406177592324"
$ws.Range("J7").Value2 = [double]"-4.918638627767542"
$ws.Range("B8").Value2 = [double]"-3.051376050698078"
$ws.Range("C8").Value2 = [double]"3.060717104569447"
$ws.Range("D8").Value2 = [double]"2.362636269252887"
$ws.Range("E8").Value2 = [double]"-2.02616217426285"
$ws.Range("F8").Value2 = [double]"-3.119675605408095"
$ws.Range("G8").Value2 = [double]"-1.311114057127361"
$ws.Range("H8").Value2 = [double]"0"
$ws.Range("I8").Value2 = [double]"1.223055507014901"
$ws.Range("J8").Value2 = [double]"-5.107079812482137"
$ws.Range("B9").Value2 = [double]"-2.955561379033963"
$ws.Range("C9").Value2 = [double]"1.902179684934643"
$ws.Range("D9").Value2 = [double]"0.6518618050328975"
$ws.Range("E9").Value2 = [double]"-2.214968684321045"
$ws.Range("F9").Value2 = [double]"-3.01932165680969"
$ws.Range("G9").Value2 = [double]"-1.433406177592324"
$ws.Range("H9").Value2 = [double]"-1.223055507014901"
$ws.Range("I9").Value2 = [double]"0"
$ws.Range("J9").Value2 = [double]"-3.682690064078619"
$ws.Range("B10").Value2 = [double]"-1.990661243152962"
$ws.Range("C10").Value2 = [double]"4.389490135705872"
$ws.Range("D10").Value2 = [double]"4.192119868523634"
$ws.Range("E10").Value2 = [double]"3.388081490838152"
$ws.Range("F10").Value2 = [double]"-1.805059303861617"
$ws.Range("G10").Value2 = [double]"4.918638627767542"
$ws.Range("H10").Value2 = [double]"5.107079812482137"
$ws.Range("I10").Value2 = [double]"3.682690064078619"
$ws.Range("J10").Value2 = [double]"0"

$ws = $wb.Worksheets.Item("Resumen")
$ws.Range("A2").Value2 = "Sieve Bootstrap"
$ws.Range("B2").Value2 = [double]"3"
$ws.Range("C2").Value2 = [double]"0"
$ws.Range("D2").Value2 = [double]"5"
$ws.Range("E2").Value2 = [double]"37.5"
$ws.Range("A3").Value2 = "LSPM"
$ws.Range("B3").Value2 = [double]"2"
$ws.Range("C3").Value2 = [double]"1"
$ws.Range("D3").Value2 = [double]"5"
$ws.Range("E3").Value2 = [double]"25"
$ws.Range("A4").Value2 = "DeepAR"
$ws.Range("B4").Value2 = [double]"1"
$ws.Range("C4").Value2 = [double]"0"
$ws.Range("D4").Value2 = [double]"7"
$ws.Range("E4").Value2 = [double]"12.5"
$ws.Range("A5").Value2 = "AV-MCPS"
$ws.Range("B5").Value2 = [double]"1"
$ws.Range("C5").Value2 = [double]"0"
$ws.Range("D5").Value2 = [double]"7"
$ws.Range("E5").Value2 = [double]"12.5"
$ws.Range("A6").Value2 = "MCPS"
$ws.Range("B6").Value2 = [double]"1"
$ws.Range("C6").Value2 = [double]"0"
$ws.Range("D6").Value2 = [double]"7"
$ws.Range("E6").Value2 = [double]"12.5"
$ws.Range("A7").Value2 = "Block Bootstrapping"
$ws.Range("B7").Value2 = [double]"0"
$ws.Range("C7").Value2 = [double]"0"
$ws.Range("D7").Value2 = [double]"8"
$ws.Range("E7").Value2 = [double]"0"
$ws.Range("A8").Value2 = "LSPMW"
$ws.Range("B8").Value2 = [double]"0"
$ws.Range("C8").Value2 = [double]"2"
$ws.Range("D8").Value2 = [double]"6"
$ws.Range("E8").Value2 = [double]"0"
$ws.Range("A9").Value2 = "AREPD"
$ws.Range("B9").Value2 = [double]"0"
$ws.Range("C9").Value2 = [double]"0"
$ws.Range("D9").Value2 = [double]"8"
$ws.Range("E9").Value2 = [double]"0"
$ws.Range("A10").Value2 = "EnCQR-LSTM"
$ws.Range("C10").Value2 = [double]"5"
$ws.Range("D10").Value2 = [double]"3"
